$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "ideal040616"
$ws.Range("B1").Value = "USD"
$ws.Range("C1").NumberFormat = "@"
$ws.Range("C1").Value = "10"
$ws.Range("C1").Style = "Normal"
$ws.Range("D1").NumberFormat = "@"
$ws.Range("D1").Value = "10"
$ws.Range("D1").Style = "Normal"
$ws.Range("E1").Value = "Visa,Worldpay"
$ws.Range("F1").NumberFormat = "@"
$ws.Range("F1").Value = "3381930175"
$ws.Range("F1").Style = "Normal"
$ws.Range("G1").NumberFormat = "@"
$ws.Range("G1").Value = "20.0"
$ws.Range("G1").Style = "Normal"
$ws.Range("H1").NumberFormat = "@"
$ws.Range("H1").Value = "20.0"
$ws.Range("H1").Style = "Normal"

$ws.Range("A2").Value = "edeal040616"
$ws.Range("B2").Value = "SAR"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "89.92"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "40"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "MasterCard,Worldpay"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "3306792453"
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "129.92"
$ws.Range("G2").Style = "Normal"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "129.92"
$ws.Range("H2").Style = "Normal"

$ws.Range("A3").Value = "edeal040616"
$ws.Range("B3").Value = "KWD"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "10.41"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "4"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "American Express,Worldpay"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "3369586302"
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "14.41"
$ws.Range("G3").Style = "Normal"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "14.41"
$ws.Range("H3").Style = "Normal"

$ws.Range("A4").Value = "idealReseller040616"
$ws.Range("B4").Value = "NZD"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "210"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "200"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "Visa,Worldpay"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "3362472968"
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "420.0"
$ws.Range("G4").Style = "Normal"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "420.0"
$ws.Range("H4").Style = "Normal"

$ws.Range("A5").Value = "iuliia.6"
$ws.Range("B5").Value = "SGD"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "38025.03"
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "175"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "MasterCard,Worldpay"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "3342276487"
$ws.Range("F5").Style = "Normal"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "38200.03"
$ws.Range("G5").Style = "Normal"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "38199.94"
$ws.Range("H5").Style = "Normal"

$ws.Range("A6").Value = "iuliia.6"
$ws.Range("B6").Value = "QAR"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "99485.96"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "500"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "American Express,Worldpay"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "3323739705"
$ws.Range("F6").Style = "Normal"
$ws.Range("G6").ClearContents()
$ws.Range("H6").ClearContents()

$ws.Range("A7").Value = "iuliia.6"
$ws.Range("B7").Value = "QAR"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "99985.93"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "500"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "American Express,Worldpay"
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "3339014984"
$ws.Range("F7").Style = "Normal"
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()

$ws.Range("A8").Value = "iuliia.6"
$ws.Range("B8").Value = "QAR"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "100486.25"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "500"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "American Express,Worldpay"
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "3337963529"
$ws.Range("F8").Style = "Normal"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "100986.25"
$ws.Range("G8").Style = "Normal"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "100986.22"
$ws.Range("H8").Style = "Normal"

$ws.Range("A9").Value = "ideal040616"
$ws.Range("B9").Value = "USD"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "20"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "10"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "Visa,Worldpay"
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "3315144745"
$ws.Range("F9").Style = "Normal"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "30.0"
$ws.Range("G9").Style = "Normal"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "30.0"
$ws.Range("H9").Style = "Normal"

$ws.Range("A10").Value = "edeal040616"
$ws.Range("B10").Value = "SAR"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "179.85"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "MasterCard,Worldpay"
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = "3329943515"
$ws.Range("F10").Style = "Normal"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "219.85"
$ws.Range("G10").Style = "Normal"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "219.85"
$ws.Range("H10").Style = "Normal"

$ws.Range("A11").Value = "edeal040616"
$ws.Range("B11").Value = "KWD"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "17.61"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "American Express,Worldpay"
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = "3369982417"
$ws.Range("F11").Style = "Normal"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "21.61"
$ws.Range("G11").Style = "Normal"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "21.61"
$ws.Range("H11").Style = "Normal"

$ws.Range("A12").Value = "idealReseller040616"
$ws.Range("B12").Value = "NZD"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "420"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "200"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "Visa,Worldpay"
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = "3325455265"
$ws.Range("F12").Style = "Normal"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "630.0"
$ws.Range("G12").Style = "Normal"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "630.0"
$ws.Range("H12").Style = "Normal"

$ws.Range("A13").Value = "iuliia.6"
$ws.Range("B13").Value = "SGD"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "38776"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "175"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "MasterCard,Worldpay"
$ws.Range("F13").NumberFormat = "@"
$ws.Range("F13").Value = "3397636624"
$ws.Range("F13").Style = "Normal"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "38951.0"
$ws.Range("G13").Style = "Normal"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "38950.92"
$ws.Range("H13").Style = "Normal"

$ws.Range("A14").Value = "iuliia.6"
$ws.Range("B14").Value = "QAR"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "101441.76"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "500"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "American Express,Worldpay"
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = "3324200528"
$ws.Range("F14").Style = "Normal"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "101941.76"
$ws.Range("G14").Style = "Normal"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "101941.72"
$ws.Range("H14").Style = "Normal"

$ws.Range("A15").Value = "ideal040616"
$ws.Range("B15").Value = "CAD"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "41.2"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "Visa,Adyen"
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = "3336774126"
$ws.Range("F15").Style = "Normal"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "55.2"
$ws.Range("G15").Style = "Normal"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "55.2"
$ws.Range("H15").Style = "Normal"

$ws.Range("A16").Value = "edeal040616"
$ws.Range("B16").Value = "AED"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "264.21"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "40"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "MasterCard,Adyen"
$ws.Range("F16").NumberFormat = "@"
$ws.Range("F16").Value = "3362518051"
$ws.Range("F16").Style = "Normal"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "304.21"
$ws.Range("G16").Style = "Normal"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "304.21"
$ws.Range("H16").Style = "Normal"

$ws.Range("A17").Value = "edeal040616"
$ws.Range("B17").Value = "KWD"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "24.89"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "American Express,Adyen"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "28.89"
$ws.Range("G17").Style = "Normal"
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "24.89"
$ws.Range("H17").Style = "Normal"

$ws.Range("A18").Value = "idealReseller040616"
$ws.Range("B18").Value = "NZD"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "630"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "200"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "Visa,Adyen"
$ws.Range("F18").NumberFormat = "@"
$ws.Range("F18").Value = "3309155197"
$ws.Range("F18").Style = "Normal"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "840.0"
$ws.Range("G18").Style = "Normal"
$ws.Range("H18").NumberFormat = "@"
$ws.Range("H18").Value = "840.0"
$ws.Range("H18").Style = "Normal"

$ws.Range("A19").Value = "iuliia.6"
$ws.Range("B19").Value = "AUD"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "38941.32"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "160"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "MasterCard,Adyen"
$ws.Range("F19").NumberFormat = "@"
$ws.Range("F19").Value = "3317205756"
$ws.Range("F19").Style = "Normal"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "39101.32"
$ws.Range("G19").Style = "Normal"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "39101.43"
$ws.Range("H19").Style = "Normal"

$ws.Range("A20").Value = "iuliia.6"
$ws.Range("B20").Value = "GBP"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "17450.87"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "125"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "American Express,Adyen"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "17575.87"
$ws.Range("G20").Style = "Normal"
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = "17450.87"
$ws.Range("H20").Style = "Normal"

$ws.Range("A21").Value = "ideal040616"
$ws.Range("B21").Value = "CAD"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "55.2"
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "Visa,Adyen"
$ws.Range("F21").NumberFormat = "@"
$ws.Range("F21").Value = "3354129786"
$ws.Range("F21").Style = "Normal"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "69.2"
$ws.Range("G21").Style = "Normal"
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "69.2"
$ws.Range("H21").Style = "Normal"

$ws.Range("A22").Value = "edeal040616"
$ws.Range("B22").Value = "AED"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "304.21"
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "MasterCard,Adyen"
$ws.Range("F22").NumberFormat = "@"
$ws.Range("F22").Value = "3325088319"
$ws.Range("F22").Style = "Normal"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "344.21"
$ws.Range("G22").Style = "Normal"
$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = "344.21"
$ws.Range("H22").Style = "Normal"

$ws.Range("A23").Value = "edeal040616"
$ws.Range("B23").Value = "KWD"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "28.16"
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "American Express,Adyen"
$ws.Range("F23").NumberFormat = "@"
$ws.Range("F23").Value = "3313456880"
$ws.Range("F23").Style = "Normal"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "32.16"
$ws.Range("G23").Style = "Normal"
$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value = "32.16"
$ws.Range("H23").Style = "Normal"

$ws.Range("A24").Value = "idealReseller040616"
$ws.Range("B24").Value = "NZD"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "840"
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "200"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "Visa,Adyen"
$ws.Range("F24").NumberFormat = "@"
$ws.Range("F24").Value = "3384860085"
$ws.Range("F24").Style = "Normal"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "1050.0"
$ws.Range("G24").Style = "Normal"
$ws.Range("H24").NumberFormat = "@"
$ws.Range("H24").Value = "1050.0"
$ws.Range("H24").Style = "Normal"

$ws.Range("A25").Value = "iuliia.6"
$ws.Range("B25").Value = "AUD"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "39101.43"
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "MasterCard,Adyen"
$ws.Range("F25").NumberFormat = "@"
$ws.Range("F25").Value = "3323988613"
$ws.Range("F25").Style = "Normal"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "39261.43"
$ws.Range("G25").Style = "Normal"
$ws.Range("H25").NumberFormat = "@"
$ws.Range("H25").Value = "39261.39"
$ws.Range("H25").Style = "Normal"

$ws.Range("A26").Value = "iuliia.6"
$ws.Range("B26").Value = "GBP"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "17522.27"
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "125"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "American Express,Adyen"
$ws.Range("F26").NumberFormat = "@"
$ws.Range("F26").Value = "3337517038"
$ws.Range("F26").Style = "Normal"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "17647.27"
$ws.Range("G26").Style = "Normal"
$ws.Range("H26").NumberFormat = "@"
$ws.Range("H26").Value = "17647.23"
$ws.Range("H26").Style = "Normal"

$ws.Range("A27").Value = "ideal040616"
$ws.Range("B27").Value = "CAD"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "69.2"
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "Visa,Global Collect"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "83.2"
$ws.Range("G27").Style = "Normal"
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H27").Value = "69.2"
$ws.Range("H27").Style = "Normal"

$ws.Range("A28").Value = "edeal040616"
$ws.Range("B28").Value = "AED"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "393.1"
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "MasterCard,Global Collect"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "397.1"
$ws.Range("G28").Style = "Normal"
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = "32.16"
$ws.Range("H28").Style = "Normal"

$ws.Range("A29").Value = "idealReseller040616"
$ws.Range("B29").Value = "NZD"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "1050"
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "200"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "Visa,Global Collect"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "1260.0"
$ws.Range("G29").Style = "Normal"
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = "1050.0"
$ws.Range("H29").Style = "Normal"

$ws.Range("A30").Value = "iuliia.6"
$ws.Range("B30").Value = "AUD"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "39541.4"
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "160"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "MasterCard,Global Collect"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "39666.4"
$ws.Range("G30").Style = "Normal"

Write-Output "done"